$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Values -----------------------------------------------------------
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "T.C (Azure)"
$ws.Range("C1").Value = "T.C (Desc.)"
$ws.Range("D1").Value = "Error"

$ws.Range("A2").Value = 12194210
$ws.Range("B2").Value = 114447

# --- Header font (bold, 14pt) ------------------------------------------
$headerFont = $ws.Range("A1:D1").Font
$headerFont.Bold = $true
$headerFont.Size = 14

# --- Row heights --------------------------------------------------------
$ws.Rows(1).RowHeight = 18.75
$ws.Rows(2).RowHeight = 15.75

# --- Thick-box border around A1:D2, thin gridlines inside ---------------
$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlInsideVertical = 11
$xlInsideHorizontal = 12

$xlContinuous = 1
$xlThin = 2
$xlMedium = -4138

$tbl = $ws.Range("A1:D2")

$tbl.Borders($xlEdgeLeft).LineStyle = $xlContinuous
$tbl.Borders($xlEdgeLeft).Weight = $xlMedium

$tbl.Borders($xlEdgeRight).LineStyle = $xlContinuous
$tbl.Borders($xlEdgeRight).Weight = $xlMedium

$tbl.Borders($xlEdgeTop).LineStyle = $xlContinuous
$tbl.Borders($xlEdgeTop).Weight = $xlMedium

$tbl.Borders($xlEdgeBottom).LineStyle = $xlContinuous
$tbl.Borders($xlEdgeBottom).Weight = $xlMedium

$tbl.Borders($xlInsideVertical).LineStyle = $xlContinuous
$tbl.Borders($xlInsideVertical).Weight = $xlThin

$tbl.Borders($xlInsideHorizontal).LineStyle = $xlContinuous
$tbl.Borders($xlInsideHorizontal).Weight = $xlThin

# --- Active cell / selection ---------------------------------------------
$ws.Range("B6").Select()
